$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (column D) and Volume(1h) (column E) cells per latest crypto snapshot ---

# Row 2
$ws.Range("D2").Value = "27.192.26"
$ws.Range("E2").Value = "  -1.03%  "

# Row 3
$ws.Range("D3").Value = "1.576.47"
$ws.Range("E3").Value = "  -0.30%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.48"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.26%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.491"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.25%  "

# Row 7
$ws.Range("E7").Value = "  -0.05%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.32"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.07%  "

# Row 9
$ws.Range("E9").Value = "  -0.20%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0590"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.09%  "

# Row 11
$ws.Range("E11").Value = "  +0.26%  "

# Row 12
$ws.Range("D12").Value = "1.799.24"
$ws.Range("E12").Value = "  -0.24%  "

# Row 13
$ws.Range("D13").Value = "1.563.66"
$ws.Range("E13").Value = "  -1.02%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.78"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.35%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.519"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.04%  "

# Row 16
$ws.Range("D16").Value = "27.187.87"
$ws.Range("E16").Value = "  -1.09%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.37"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.20%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.40"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.06%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "214.69"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.09%  "

# Row 20
$ws.Range("E20").Value = "  -0.77%  "

# Row 21
$ws.Range("E21").Value = "  +0.06%  "

# Row 22
$ws.Range("E22").Value = "  -0.22%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.51"
$ws.Range("D23").ClearFormats()

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.02"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.81%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.58"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.50%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.71"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.02%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.98"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.52%  "

# Row 28
$ws.Range("E28").Value = "  -0.01%  "

# Row 29
$ws.Range("E29").Value = "  -1.10%  "

# Row 30
$ws.Range("E30").Value = "  -3.27%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0464"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.68%  "

# Row 32
$ws.Range("E32").Value = "  -1.00%  "

# Row 33
$ws.Range("D33").Value = "1.400.24"
$ws.Range("E33").Value = "  +2.18%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.92"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.03%  "

# Row 35
$ws.Range("E35").Value = "  +1.36%  "

# Row 36
$ws.Range("E36").Value = "  -1.16%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.940"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.85%  "

# Row 38
$ws.Range("E38").Value = "  -1.92%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.818"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.87%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.516"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.90%  "

# Row 41
$ws.Range("E41").Value = "  +0.04%  "

# Row 42
$ws.Range("E42").Value = "  +3.46%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.84"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.04%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.46"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.48%  "

# Row 45
$ws.Range("E45").Value = "  +1.20%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.87"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.86%  "

# Row 47
$ws.Range("D47").Value = "1.711.74"
$ws.Range("E47").Value = "  -0.20%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.83"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.28%  "

# Row 49
$ws.Range("D49").Value = "0.0₇0979"
$ws.Range("E49").Value = "  -1.41%  "

# Row 50
$ws.Range("E50").Value = "  -0.61%  "

# Row 51
$ws.Range("E51").Value = "  +0.09%  "
